# Cập nhật cho bản phát hành v2.2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1-3 got slightly taller (18.75pt -> 19.5pt)
$ws.Range("A1:B3").RowHeight = 19.5

# Row 2: "NMCD 6 loại" / "SLNMCD - ..." renamed to the new UTN HQ report
$ws.Range("A2").Value = "UTNHQTT"
$ws.Range("B2").Value = "SLUTNHQ90190 - Sản lượng UTN HQ truyền thống 90g - 190g"

# Row 3: old "DHTC" entries cleared out, leaving a text quote-prefix (empty, forced text)
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "'"
